$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "MM260" origin marker in column B for every data row (3-98),
# mirroring column A's label for each row.
for ($r = 3; $r -le 98; $r++) {
    $ws.Cells.Item($r, 2).Value = "MM260"
}

# Update the sheet view so it reflects the new scroll position / selection
# that resulted from editing the sheet (selecting the newly added B column).
$ws.Range("B3:B98").Select()
$excel.ActiveWindow.ScrollRow = 77
